# Auto-generated: updates hardcoded market-price snapshot values
# across the 8 Leve-profit worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Values are plain numbers (no formulas) refreshed by a scheduled data-pull run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 30657
$ws.Range("J3").Value = 30657
$ws.Range("L3").Value = 30657
$ws.Range("N3").Value = -30885
$ws.Range("H4").Value = 400
$ws.Range("I4").Value = 400
$ws.Range("K4").Value = 400
$ws.Range("M4").Value = -286
$ws.Range("H102").Value = 30657
$ws.Range("J102").Value = 30657
$ws.Range("L102").Value = 30657
$ws.Range("N102").Value = -37147
$ws.Range("H107").Value = 335.95
$ws.Range("I107").Value = 279.6316
$ws.Range("K107").Value = 279.6316
$ws.Range("M107").Value = 1640.3684
$ws.Range("H132").Value = 3540.611
$ws.Range("I132").Value = 3683.2292
$ws.Range("K132").Value = 11049.6876
$ws.Range("M132").Value = -8519.687600000001
$ws.Range("H137").Value = 1391962.4
$ws.Range("I137").Value = 3573363.2
$ws.Range("J137").Value = 3798.1365
$ws.Range("K137").Value = 10720089.6
$ws.Range("L137").Value = 11394.4095
$ws.Range("M137").Value = -10717539.6
$ws.Range("N137").Value = -16494.4095

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4444.6343
$ws.Range("I32").Value = 3963.2942
$ws.Range("K32").Value = 3963.2942
$ws.Range("M32").Value = -3676.2942
$ws.Range("H45").Value = 41262.184
$ws.Range("I45").Value = 61404
$ws.Range("K45").Value = 61404
$ws.Range("M45").Value = -61027
$ws.Range("H74").Value = 233864.55
$ws.Range("I74").Value = 618191.7
$ws.Range("K74").Value = 618191.7
$ws.Range("M74").Value = -617317.7
$ws.Range("H77").Value = 233864.55
$ws.Range("I77").Value = 618191.7
$ws.Range("K77").Value = 3090958.5
$ws.Range("M77").Value = -3086590.5
$ws.Range("H97").Value = 1176.8667
$ws.Range("I97").Value = 1076.7727
$ws.Range("J97").Value = 1452.125
$ws.Range("K97").Value = 1076.7727
$ws.Range("L97").Value = 1452.125
$ws.Range("M97").Value = -580.7727
$ws.Range("N97").Value = -2444.125
$ws.Range("H104").Value = 65000
$ws.Range("J104").Value = 65000
$ws.Range("L104").Value = 65000
$ws.Range("N104").Value = -71988
$ws.Range("H122").Value = 19998.334
$ws.Range("J122").Value = 19998.5
$ws.Range("L122").Value = 59995.5
$ws.Range("N122").Value = -64895.5
$ws.Range("H132").Value = 2132.3704
$ws.Range("I132").Value = 1033.8823
$ws.Range("J132").Value = 3999.8
$ws.Range("K132").Value = 3101.6469
$ws.Range("L132").Value = 11999.4
$ws.Range("M132").Value = -571.6468999999997
$ws.Range("N132").Value = -17059.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2088.6365
$ws.Range("I86").Value = 1820
$ws.Range("J86").Value = 2411
$ws.Range("K86").Value = 1820
$ws.Range("L86").Value = 2411
$ws.Range("M86").Value = -697
$ws.Range("N86").Value = -4657
$ws.Range("H89").Value = 2088.6365
$ws.Range("I89").Value = 1820
$ws.Range("J89").Value = 2411
$ws.Range("K89").Value = 9100
$ws.Range("L89").Value = 12055
$ws.Range("M89").Value = -3484
$ws.Range("N89").Value = -23287
$ws.Range("H105").Value = 17335488
$ws.Range("I105").Value = 1001928.3
$ws.Range("J105").Value = 50002610
$ws.Range("K105").Value = 1001928.3
$ws.Range("L105").Value = 50002610
$ws.Range("M105").Value = -1000181.3
$ws.Range("N105").Value = -50006104
$ws.Range("H125").Value = 103665
$ws.Range("J125").Value = 103665
$ws.Range("L125").Value = 103665
$ws.Range("N125").Value = -113505

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1501.8928
$ws.Range("I16").Value = 1323.5454
$ws.Range("K16").Value = 1323.5454
$ws.Range("M16").Value = -1036.5454
$ws.Range("H31").Value = 5303.684
$ws.Range("I31").Value = 3480.8333
$ws.Range("J31").Value = 8428.571
$ws.Range("K31").Value = 3480.8333
$ws.Range("L31").Value = 8428.571
$ws.Range("M31").Value = -3185.8333
$ws.Range("N31").Value = -9018.571
$ws.Range("H34").Value = 5303.684
$ws.Range("I34").Value = 3480.8333
$ws.Range("J34").Value = 8428.571
$ws.Range("K34").Value = 3480.8333
$ws.Range("L34").Value = 8428.571
$ws.Range("M34").Value = -3278.8333
$ws.Range("N34").Value = -8832.571
$ws.Range("H94").Value = 1960.9375
$ws.Range("J94").Value = 2044.5
$ws.Range("L94").Value = 2044.5
$ws.Range("N94").Value = -2946.5
$ws.Range("H113").Value = 1501.8928
$ws.Range("I113").Value = 1323.5454
$ws.Range("K113").Value = 1323.5454
$ws.Range("M113").Value = 846.4546
$ws.Range("H122").Value = 4487.25
$ws.Range("I122").Value = 3768.0625
$ws.Range("K122").Value = 11304.1875
$ws.Range("M122").Value = -8854.1875
$ws.Range("H132").Value = 2004.9375
$ws.Range("I132").Value = 1673.4166
$ws.Range("K132").Value = 5020.2498
$ws.Range("M132").Value = -2490.2498

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 116618.664
$ws.Range("J25").Value = 131071
$ws.Range("L25").Value = 393213
$ws.Range("N25").Value = -393551
$ws.Range("H30").Value = 116618.664
$ws.Range("J30").Value = 131071
$ws.Range("L30").Value = 393213
$ws.Range("N30").Value = -393417
$ws.Range("H122").Value = 669.5
$ws.Range("I122").Value = 823.5714
$ws.Range("J122").Value = 586.53845
$ws.Range("K122").Value = 7412.1426
$ws.Range("L122").Value = 5278.84605
$ws.Range("M122").Value = -4962.1426
$ws.Range("N122").Value = -10178.84605
$ws.Range("H129").Value = 555000
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 555000
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 1665000
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -1675000

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 117381.22
$ws.Range("J70").Value = 7807.5454
$ws.Range("L70").Value = 7807.5454
$ws.Range("N70").Value = -8347.545399999999
$ws.Range("H73").Value = 117381.22
$ws.Range("J73").Value = 7807.5454
$ws.Range("L73").Value = 7807.5454
$ws.Range("N73").Value = -9679.545399999999
$ws.Range("H102").Value = 1339.25
$ws.Range("I102").Value = 1155.5714
$ws.Range("K102").Value = 1155.5714
$ws.Range("M102").Value = 466.4286
$ws.Range("H122").Value = 4367.7095
$ws.Range("I122").Value = 3108.8462
$ws.Range("J122").Value = 5276.8887
$ws.Range("K122").Value = 9326.5386
$ws.Range("L122").Value = 15830.6661
$ws.Range("M122").Value = -6876.5386
$ws.Range("N122").Value = -20730.6661
$ws.Range("H132").Value = 1051.7693
$ws.Range("I132").Value = 1051.7693
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3155.3079
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -625.3078999999998
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4178.4136
$ws.Range("I40").Value = 4246.107
$ws.Range("J40").Value = 4115.2334
$ws.Range("K40").Value = 4246.107
$ws.Range("L40").Value = 4115.2334
$ws.Range("M40").Value = -4110.107
$ws.Range("N40").Value = -4387.2334
$ws.Range("H61").Value = 14522.111
$ws.Range("I61").Value = 3065.6667
$ws.Range("K61").Value = 3065.6667
$ws.Range("M61").Value = -2863.6667
$ws.Range("H113").Value = 14522.111
$ws.Range("I113").Value = 3065.6667
$ws.Range("K113").Value = 3065.6667
$ws.Range("M113").Value = -895.6667000000002
$ws.Range("H122").Value = 13199.2
$ws.Range("I122").Value = 12999.25
$ws.Range("J122").Value = 13999
$ws.Range("K122").Value = 38997.75
$ws.Range("L122").Value = 41997
$ws.Range("M122").Value = -36547.75
$ws.Range("N122").Value = -46897
$ws.Range("H127").Value = 90874.75
$ws.Range("J127").Value = 90874.75
$ws.Range("L127").Value = 90874.75
$ws.Range("N127").Value = -100794.75
$ws.Range("H132").Value = 5386.6665
$ws.Range("I132").Value = 2880
$ws.Range("J132").Value = 5700
$ws.Range("K132").Value = 8640
$ws.Range("L132").Value = 17100
$ws.Range("M132").Value = -6110
$ws.Range("N132").Value = -22160
$ws.Range("H136").Value = 4875.9287
$ws.Range("I136").Value = 5473.778
$ws.Range("K136").Value = 16421.334
$ws.Range("M136").Value = -13871.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 250500
$ws.Range("J68").Value = 250500
$ws.Range("L68").Value = 250500
$ws.Range("N68").Value = -252122
$ws.Range("H71").Value = 250500
$ws.Range("J71").Value = 250500
$ws.Range("L71").Value = 751500
$ws.Range("N71").Value = -759612
$ws.Range("H107").Value = 1028.75
$ws.Range("I107").Value = 1028.75
$ws.Range("K107").Value = 3086.25
$ws.Range("M107").Value = -1166.25
$ws.Range("H122").Value = 6580002.5
$ws.Range("I122").Value = 793.1724
$ws.Range("K122").Value = 2379.5172
$ws.Range("M122").Value = 70.48279999999977
$ws.Range("H132").Value = 3264.2632
$ws.Range("I132").Value = 3220.0625
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 9660.1875
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -7130.1875
$ws.Range("N132").Value = -15560
$ws.Range("H136").Value = 71429656
$ws.Range("I136").Value = 111112000
$ws.Range("K136").Value = 333336000
$ws.Range("M136").Value = -333333450

